# This edit inserts one new data row for "Cereza" (variety Lapins) right
# before the existing row 773 in the Mercado Mayorista Lo Valledor de
# Santiago weekly price sheet. Inserting a full row pushes every
# subsequent row (old 773..805) down by one (new 774..806), which is
# exactly what the target diff shows (old row 773 data now lives at 774,
# etc., and a brand new row appears at 773; dimension grows to A1:T806).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 773, shifting rows
# 773:805 down to 774:806.
$ws.Rows(773).Insert()

# Populate the newly inserted row 773 with the new record.
$ws.Range("A773").Value = 6
$ws.Range("B773").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C773").Value = "Metropolitana"
$ws.Range("D773").Value = 44595
$ws.Range("E773").Value = 13
$ws.Range("F773").Value = "Fruta"
$ws.Range("G773").Value = 100103
$ws.Range("H773").Value = "Frutos de hueso (carozo)"
$ws.Range("I773").Value = 100103001
$ws.Range("J773").Value = "Cereza"
$ws.Range("K773").Value = "Lapins"
$ws.Range("L773").Value = "Primera"
$ws.Range("M773").Value = 250
$ws.Range("N773").Value = 8000
$ws.Range("O773").Value = 8000
$ws.Range("P773").Value = 8000
$ws.Range("Q773").Value = "`$/bandeja 10 kilos"
$ws.Range("R773").Value = "Región de O'Higgins"
$ws.Range("S773").Value = 800
$ws.Range("T773").Value = 10
